$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$reqFerment = "LOT2028 -  Tecnologia de Processos Fermentativos  (Requisito fraco)`n"
$reqBebidas = "LOT2052 -  Tecnologia de Bebidas Experimental  (Indicação de Conjunto)`n"

# Swap the two requirement rows: the "Bebidas Experimental" entry now comes
# first (row 24), followed by the "Processos Fermentativos" entry (row 25).
$ws.Range("B24").Value = $reqBebidas
$ws.Range("C24").Value = $reqBebidas
$ws.Range("B25").Value = $reqFerment
$ws.Range("C25").Value = $reqFerment
